# Update countries & provincias Spain
# Applies the data refresh that was scraped at 20:44 (previously 19:27):
#  - Updates numeric stats for a handful of countries
#  - Marruecos overtakes Belgica -> rows 37/38 swap contents
#  - Guinea Ecuatorial overtakes Congo -> rows 122/123 swap contents
#  - Montserrat overtakes Islas Malvinas -> rows 214/215 swap contents

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 20:44"

# Helper data: row -> (Country, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes)
$rows = @(
    @{ Row = 4;   Pais = "Estados Unidos";     B = 7060345; C = 14129; D = 4314339; E = 2541081; F = 0; G = 423;  H = 204925 }
    @{ Row = 5;   Pais = "India";              B = 5640441; C = 80336; D = 4581746; E = 968674;  F = 0; G = 1056; H = 90021 }
    @{ Row = 14;  Pais = "Francia";            B = 468069;  C = 10008; D = 93538;   E = 343115;  F = 0; G = 78;   H = 31416 }
    @{ Row = 25;  Pais = "Alemania";           B = 276829;  C = 1278;  D = 246300;  E = 21040;   F = 0; G = 8;    H = 9489 }
    @{ Row = 37;  Pais = "Marruecos";          B = 105346;  C = 2227;  D = 85883;   E = 17574;   F = 0; G = 34;   H = 1889 }
    @{ Row = 38;  Pais = "Belgica";            B = 103392;  C = 1097;  D = 18977;   E = 74465;   F = 0; G = 2;    H = 9950 }
    @{ Row = 51;  Pais = "Etiopia";            B = 70422;   C = 713;   D = 28991;   E = 40304;   F = 0; G = 19;   H = 1127 }
    @{ Row = 114; Pais = "Malaui";             B = 5739;    C = 6;     D = 4065;    E = 1495;    F = 0; G = 0;    H = 179 }
    @{ Row = 117; Pais = "Suazilandia";        B = 5307;    C = 25;    D = 4672;    E = 529;     F = 0; G = 2;    H = 106 }
    @{ Row = 122; Pais = "Guinea Ecuatorial";  B = 5018;    C = 16;    D = 4509;    E = 426;     F = 0; G = 0;    H = 83 }
    @{ Row = 123; Pais = "Congo";              B = 5002;    C = 0;     D = 3887;    E = 1026;    F = 0; G = 0;    H = 89 }
    @{ Row = 143; Pais = "Mali";               B = 3030;    C = 6;     D = 2380;    E = 521;     F = 0; G = 1;    H = 129 }
    @{ Row = 214; Pais = "Montserrat";         B = 13;      C = 0;     D = 12;      E = 0;       F = 0; G = 0;    H = 1 }
    @{ Row = 215; Pais = "Islas Malvinas";     B = 13;      C = 0;     D = 13;      E = 0;       F = 0; G = 0;    H = 0 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Pais
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
    $ws.Range("H$r").Value = $item.H
}
